$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODAY")

# Update product description text (spf50+ -> spf30+)
$ws.Range("E32").Value = "Nivea Sun® Kids Spray TRIGGER spf30+ 300ml"

# Update price column (G) values for the affected rows
$priceUpdates = @{
    2  = 13.9
    3  = 13.9
    7  = 14.95
    8  = 15.2
    9  = 15.45
    11 = 15.98
    22 = 13.9
    23 = 13.55
    26 = 9.9
    27 = 16.2
    29 = 9.9
    32 = 13.4
    33 = 14.2
    34 = 15.2
    39 = 9.9
    42 = 17.9
    46 = 14.8
    47 = 14.95
    49 = 14.9
    53 = 16.4
    60 = 10.9
    76 = 6.95
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Range("G$row").Value = $priceUpdates[$row]
}
